# Generate Report for Handoff
# Updates the localization-status workbook with a new handoff GUID,
# new handoff file hashes, and refreshed handoff timestamps.

$wb = $excel.ActiveWorkbook

$oldGuid = "22edef11-9c98-4357-82d0-65c5042a4433"
$newGuid = "9dc0f76c-ace8-426d-aa45-2b4a05669e58"

$oldHash = "5f113b0a7acbe46d3d86df6745fec51a556a914e"
$newHash = "ed3090ee772d2a96b2cd0c206473cc273acf8c6e"

$newHandoffDate = "2016-03-21 12:52:26"
$newZhHandoffDate = "2016-03-21 12:52:23"

$newMdName = "$newGuid.md"
$newZhXlfName = "$newGuid.$newHash.zh-cn.xlf"
$newDeXlfName = "$newGuid.$newHash.de-de.xlf"

# Hyperlink target (Address) URLs are unchanged by this edit - only the
# cell text / hyperlink display text changes. Keep the original addresses.
$mdAddress = "https://github.com/OpenLocalizationTest/oltest/blob/a245d120449e29b462ddaa7b8500d77cf288542f/e2e/$oldGuid.md"
$zhXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7bf0a5b55c781eb0fc55e163efdd5f33bc4f1070/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldGuid.$oldHash.zh-cn.xlf"
$deXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8ff523fb1aaf9cbe9a9b805c04ccafdd6b2a4803/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldGuid.$oldHash.de-de.xlf"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Hyperlinks.Delete()
$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Range("D2").Value = $newHandoffDate
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdAddress, "", "", $newMdName)

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Hyperlinks.Delete()
$wsZh.Range("D2").Hyperlinks.Delete()
$wsZh.Range("A2").Value = $newMdName
$wsZh.Range("D2").Value = $newZhXlfName
$wsZh.Range("E2").Value = $newZhHandoffDate
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $mdAddress, "", "", $newMdName)
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), $zhXlfAddress, "", "", $newZhXlfName)

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Hyperlinks.Delete()
$wsDe.Range("D2").Hyperlinks.Delete()
$wsDe.Range("A2").Value = $newMdName
$wsDe.Range("D2").Value = $newDeXlfName
$wsDe.Range("E2").Value = $newHandoffDate
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $mdAddress, "", "", $newMdName)
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), $deXlfAddress, "", "", $newDeXlfName)
